# Регион.xlsx - apply the recorded edit:
#  - F3 row gains a new value in column G (740)
#  - Row 5 gains values in D5 (400) and G5 (800), and C5/E5/F5 keep their
#    existing values (300 / 1 / 619)
#  - L2's cached AVERAGE(D2:D19) result updates automatically (recalculated
#    below) once D5 is populated
#  - the active selection moves to G11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new cell G3
$ws.Range("G3").Value = 740

# Row 5: new cells D5 and G5; C5/E5/F5 already existed and keep their values
$ws.Range("C5").Value = 300
$ws.Range("D5").Value = 400
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 619
$ws.Range("G5").Value = 800

# Recalculate so L2 (=AVERAGE(D2:D19)) picks up the new D5 value (380 -> 385)
$excel.Calculate()

# Update the active cell / selection to G11
$ws.Range("G11").Select()
